$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Target cluster) changes from "ECs" to "MuSCs" for rows 2-4
$ws.Range("D2").Value = "MuSCs"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("D4").Value = "MuSCs"

# Row 2 updates
$ws.Range("G2").Value = 16.782487
$ws.Range("H2").Value = 50.347461
$ws.Range("I2").Value = 0.9252099721531751
$ws.Range("J2").Value = 0.9252099721531751
$ws.Range("M2").Value = 0.01393633333333333
$ws.Range("N2").Value = 0.041809
$ws.Range("Q2").Value = 0.2338863329943333
$ws.Range("R2").Value = 2.104976996949
$ws.Range("S2").Value = 0.9252099721531751
$ws.Range("T2").Value = 0.9252099721531751

# Row 3 updates
$ws.Range("I3").Value = 0.0601477084462148
$ws.Range("J3").Value = 0.0601477084462148
$ws.Range("M3").Value = 0.01393633333333333
$ws.Range("N3").Value = 0.041809
$ws.Range("Q3").Value = 0.01520490201133333
$ws.Range("R3").Value = 0.136844118102
$ws.Range("S3").Value = 0.0601477084462148
$ws.Range("T3").Value = 0.0601477084462148

# Row 4 updates
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2655986666666667
$ws.Range("H4").Value = 0.7967960000000001
$ws.Range("I4").Value = 0.01464231940061012
$ws.Range("J4").Value = 0.01464231940061012
$ws.Range("M4").Value = 0.01393633333333333
$ws.Range("N4").Value = 0.041809
$ws.Range("Q4").Value = 0.003701471551555556
$ws.Range("R4").Value = 0.033313243964
$ws.Range("S4").Value = 0.01464231940061012
$ws.Range("T4").Value = 0.01464231940061012
